$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3162.5122
$ws.Range("I64").Value = 2985.9033
$ws.Range("J64").Value = 3710
$ws.Range("K64").Value = 2985.9033
$ws.Range("L64").Value = 3710
$ws.Range("M64").Value = -2737.9033
$ws.Range("N64").Value = -4206

$ws.Range("H67").Value = 3162.5122
$ws.Range("I67").Value = 2985.9033
$ws.Range("J67").Value = 3710
$ws.Range("K67").Value = 2985.9033
$ws.Range("L67").Value = 3710
$ws.Range("M67").Value = -2127.9033
$ws.Range("N67").Value = -5426

$ws.Range("H74").Value = 3865.3044
$ws.Range("I74").Value = 3703
$ws.Range("K74").Value = 3703
$ws.Range("M74").Value = -2767

$ws.Range("H76").Value = 10180
$ws.Range("J76").Value = 4125
$ws.Range("L76").Value = 4125
$ws.Range("N76").Value = -4755

$ws.Range("H77").Value = 3865.3044
$ws.Range("I77").Value = 3703
$ws.Range("K77").Value = 18515
$ws.Range("M77").Value = -13835

$ws.Range("H79").Value = 10180
$ws.Range("J79").Value = 4125
$ws.Range("L79").Value = 4125
$ws.Range("N79").Value = -6309

$ws.Range("H135").Value = 26870.41
$ws.Range("I135").Value = 42848.832
$ws.Range("J135").Value = 1304.9333
$ws.Range("K135").Value = 385639.488
$ws.Range("L135").Value = 11744.3997
$ws.Range("M135").Value = -383104.488
$ws.Range("N135").Value = -16814.3997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2135.6667
$ws.Range("I61").Value = 1908.2778
$ws.Range("J61").Value = 2590.4443
$ws.Range("K61").Value = 1908.2778
$ws.Range("L61").Value = 2590.4443
$ws.Range("M61").Value = -1696.2778
$ws.Range("N61").Value = -3014.4443

$ws.Range("H88").Value = 2045.8334
$ws.Range("I88").Value = 1275
$ws.Range("J88").Value = 2431.25
$ws.Range("K88").Value = 1275
$ws.Range("L88").Value = 2431.25
$ws.Range("M88").Value = -869
$ws.Range("N88").Value = -3243.25

$ws.Range("H91").Value = 2045.8334
$ws.Range("I91").Value = 1275
$ws.Range("J91").Value = 2431.25
$ws.Range("K91").Value = 1275
$ws.Range("L91").Value = 2431.25
$ws.Range("M91").Value = 129
$ws.Range("N91").Value = -5239.25

$ws.Range("H97").Value = 3832.7273
$ws.Range("I97").Value = 2757.5
$ws.Range("K97").Value = 2757.5
$ws.Range("M97").Value = -2261.5

$ws.Range("H136").Value = 2135.6667
$ws.Range("I136").Value = 1908.2778
$ws.Range("J136").Value = 2590.4443
$ws.Range("K136").Value = 5724.8334
$ws.Range("L136").Value = 7771.3329
$ws.Range("M136").Value = -3174.8334
$ws.Range("N136").Value = -12871.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 3093.3333
$ws.Range("I113").Value = 3093.3333
$ws.Range("K113").Value = 3093.3333
$ws.Range("M113").Value = -923.3332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1632.6923
$ws.Range("I58").Value = 1491.3684
$ws.Range("J58").Value = 2016.2858
$ws.Range("K58").Value = 1491.3684
$ws.Range("L58").Value = 2016.2858
$ws.Range("M58").Value = -1288.3684
$ws.Range("N58").Value = -2422.2858

$ws.Range("H74").Value = 34479.43
$ws.Range("J74").Value = 34479.43
$ws.Range("L74").Value = 34479.43
$ws.Range("N74").Value = -36227.43

$ws.Range("H77").Value = 34479.43
$ws.Range("J77").Value = 34479.43
$ws.Range("L77").Value = 103438.29
$ws.Range("N77").Value = -112174.29

$ws.Range("H136").Value = 1632.6923
$ws.Range("I136").Value = 1491.3684
$ws.Range("J136").Value = 2016.2858
$ws.Range("K136").Value = 4474.1052
$ws.Range("L136").Value = 6048.857400000001
$ws.Range("M136").Value = -1924.1052
$ws.Range("N136").Value = -11148.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 603.9394
$ws.Range("I5").Value = 334.58334
$ws.Range("J5").Value = 1322.2222
$ws.Range("K5").Value = 1003.75002
$ws.Range("L5").Value = 3966.6666
$ws.Range("M5").Value = -891.7500200000001
$ws.Range("N5").Value = -4190.6666

$ws.Range("H113").Value = 598.5769
$ws.Range("I113").Value = 591.125
$ws.Range("J113").Value = 601.8889
$ws.Range("K113").Value = 1773.375
$ws.Range("L113").Value = 1805.6667
$ws.Range("M113").Value = 396.625
$ws.Range("N113").Value = -6145.6667

$ws.Range("H135").Value = 603.9394
$ws.Range("I135").Value = 334.58334
$ws.Range("J135").Value = 1322.2222
$ws.Range("K135").Value = 3011.25006
$ws.Range("L135").Value = 11899.9998
$ws.Range("M135").Value = -476.2500600000003
$ws.Range("N135").Value = -16969.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3208.1
$ws.Range("I80").Value = 2927.7778
$ws.Range("J80").Value = 3437.4546
$ws.Range("K80").Value = 2927.7778
$ws.Range("L80").Value = 3437.4546
$ws.Range("M80").Value = -1929.7778
$ws.Range("N80").Value = -5433.4546

$ws.Range("H83").Value = 3208.1
$ws.Range("I83").Value = 2927.7778
$ws.Range("J83").Value = 3437.4546
$ws.Range("K83").Value = 14638.889
$ws.Range("L83").Value = 17187.273
$ws.Range("M83").Value = -9646.888999999999
$ws.Range("N83").Value = -27171.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 818.27026
$ws.Range("I22").Value = 648.4286
$ws.Range("J22").Value = 921.65216
$ws.Range("K22").Value = 648.4286
$ws.Range("L22").Value = 921.65216
$ws.Range("M22").Value = -353.4286
$ws.Range("N22").Value = -1511.65216

$ws.Range("H27").Value = 818.27026
$ws.Range("I27").Value = 648.4286
$ws.Range("J27").Value = 921.65216
$ws.Range("K27").Value = 648.4286
$ws.Range("L27").Value = 921.65216
$ws.Range("M27").Value = -541.4286
$ws.Range("N27").Value = -1135.65216

$ws.Range("H61").Value = 1918.8948
$ws.Range("I61").Value = 1582.1818
$ws.Range("J61").Value = 2381.875
$ws.Range("K61").Value = 1582.1818
$ws.Range("L61").Value = 2381.875
$ws.Range("M61").Value = -1380.1818
$ws.Range("N61").Value = -2785.875

$ws.Range("H113").Value = 1918.8948
$ws.Range("I113").Value = 1582.1818
$ws.Range("J113").Value = 2381.875
$ws.Range("K113").Value = 1582.1818
$ws.Range("L113").Value = 2381.875
$ws.Range("M113").Value = 587.8181999999999
$ws.Range("N113").Value = -6721.875

$ws.Range("H132").Value = 2515.6
$ws.Range("I132").Value = 2160.1538
$ws.Range("J132").Value = 2900.6667
$ws.Range("K132").Value = 6480.4614
$ws.Range("L132").Value = 8702.000100000001
$ws.Range("M132").Value = -3950.4614
$ws.Range("N132").Value = -13762.0001
